$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 311
$ws.Range("I12").Value = 311
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 311
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -141
$ws.Range("N12").ClearContents()
$ws.Range("H33").Value = 5292.643
$ws.Range("I33").Value = 6345.8887
$ws.Range("J33").Value = 3396.8
$ws.Range("K33").Value = 6345.8887
$ws.Range("L33").Value = 3396.8
$ws.Range("M33").Value = -6116.8887
$ws.Range("N33").Value = -3854.8
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H74").Value = 14199.8
$ws.Range("I74").Value = 14199.8
$ws.Range("K74").Value = 14199.8
$ws.Range("M74").Value = -13263.8
$ws.Range("H77").Value = 14199.8
$ws.Range("I77").Value = 14199.8
$ws.Range("K77").Value = 70999
$ws.Range("M77").Value = -66319
$ws.Range("H98").Value = 4899.3335
$ws.Range("I98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("H112").Value = 2256.3333
$ws.Range("J112").Value = 2564
$ws.Range("L112").Value = 7692
$ws.Range("N112").Value = -9908
$ws.Range("H116").Value = 5145.4546
$ws.Range("I116").Value = 5500.25
$ws.Range("J116").Value = 4199.3335
$ws.Range("K116").Value = 5500.25
$ws.Range("L116").Value = 4199.3335
$ws.Range("M116").Value = -2058.25
$ws.Range("N116").Value = -11083.3335
$ws.Range("H122").Value = 4899.3335
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H130").Value = 120000
$ws.Range("J130").Value = 120000
$ws.Range("L130").Value = 120000
$ws.Range("N130").Value = -130040
$ws.Range("H133").Value = 120780
$ws.Range("J133").Value = 120780
$ws.Range("L133").Value = 120780
$ws.Range("N133").Value = -130900
$ws.Range("H136").Value = 89999
$ws.Range("J136").Value = 89999
$ws.Range("L136").Value = 89999
$ws.Range("N136").Value = -100199
$ws.Range("H138").Value = 4109.2075
$ws.Range("I138").Value = 2297.4
$ws.Range("J138").Value = 4530.558
$ws.Range("K138").Value = 6892.200000000001
$ws.Range("L138").Value = 13591.674
$ws.Range("M138").Value = -1752.200000000001
$ws.Range("N138").Value = -23871.674

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 14653
$ws.Range("I46").Value = 10733.25
$ws.Range("J46").Value = 17788.8
$ws.Range("K46").Value = 10733.25
$ws.Range("L46").Value = 17788.8
$ws.Range("M46").Value = -10414.25
$ws.Range("N46").Value = -18426.8
$ws.Range("H61").Value = 5510.8184
$ws.Range("I61").Value = 5510.8184
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 5510.8184
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -5298.8184
$ws.Range("N61").ClearContents()
$ws.Range("H74").Value = 2165.7727
$ws.Range("I74").Value = 2197.524
$ws.Range("K74").Value = 2197.524
$ws.Range("M74").Value = -1323.524
$ws.Range("H77").Value = 2165.7727
$ws.Range("I77").Value = 2197.524
$ws.Range("K77").Value = 10987.62
$ws.Range("M77").Value = -6619.619999999999
$ws.Range("H132").Value = 6471.75
$ws.Range("I132").Value = 6471.75
$ws.Range("K132").Value = 19415.25
$ws.Range("M132").Value = -16885.25
$ws.Range("H136").Value = 5510.8184
$ws.Range("I136").Value = 5510.8184
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 16532.4552
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -13982.4552
$ws.Range("N136").ClearContents()
$ws.Range("H138").Value = 88822.75
$ws.Range("J138").Value = 88822.75
$ws.Range("L138").Value = 88822.75
$ws.Range("N138").Value = -99102.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 1502.1666
$ws.Range("J12").Value = 1974.5
$ws.Range("L12").Value = 1974.5
$ws.Range("N12").Value = -2310.5
$ws.Range("H63").Value = 55135
$ws.Range("J63").Value = 55135
$ws.Range("L63").Value = 55135
$ws.Range("N63").Value = -56507
$ws.Range("H66").Value = 55135
$ws.Range("J66").Value = 55135
$ws.Range("L66").Value = 165405
$ws.Range("N66").Value = -172269
$ws.Range("H99").Value = 2800.2104
$ws.Range("I99").Value = 2107.2
$ws.Range("K99").Value = 2107.2
$ws.Range("M99").Value = -609.1999999999998
$ws.Range("H100").Value = 50000
$ws.Range("J100").Value = 50000
$ws.Range("L100").Value = 50000
$ws.Range("N100").Value = -52164
$ws.Range("H107").Value = 2533.9333
$ws.Range("I107").Value = 2203.875
$ws.Range("J107").Value = 2911.1428
$ws.Range("K107").Value = 2203.875
$ws.Range("L107").Value = 2911.1428
$ws.Range("M107").Value = -283.875
$ws.Range("N107").Value = -6751.1428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2132.262
$ws.Range("I58").Value = 2067.647
$ws.Range("K58").Value = 2067.647
$ws.Range("M58").Value = -1864.647
$ws.Range("H86").Value = 23833.5
$ws.Range("I86").Value = 38167
$ws.Range("J86").Value = 9500
$ws.Range("K86").Value = 38167
$ws.Range("L86").Value = 9500
$ws.Range("M86").Value = -37044
$ws.Range("N86").Value = -11746
$ws.Range("H89").Value = 23833.5
$ws.Range("I89").Value = 38167
$ws.Range("J89").Value = 9500
$ws.Range("K89").Value = 190835
$ws.Range("L89").Value = 47500
$ws.Range("M89").Value = -185219
$ws.Range("N89").Value = -58732
$ws.Range("H105").Value = 1518.7
$ws.Range("I105").Value = 1915.2
$ws.Range("J105").Value = 1122.2
$ws.Range("K105").Value = 1915.2
$ws.Range("L105").Value = 1122.2
$ws.Range("M105").Value = -168.2
$ws.Range("N105").Value = -4616.2
$ws.Range("H135").Value = 89999
$ws.Range("J135").Value = 89999
$ws.Range("L135").Value = 89999
$ws.Range("N135").Value = -100139
$ws.Range("H136").Value = 2132.262
$ws.Range("I136").Value = 2067.647
$ws.Range("K136").Value = 6202.941
$ws.Range("M136").Value = -3652.941

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1536.6154
$ws.Range("I11").Value = 608.55554
$ws.Range("K11").Value = 1825.66662
$ws.Range("M11").Value = -1685.66662
$ws.Range("H29").Value = 14482.143
$ws.Range("I29").Value = 33388.668
$ws.Range("J29").Value = 302.25
$ws.Range("K29").Value = 100166.004
$ws.Range("L29").Value = 906.75
$ws.Range("M29").Value = -99889.00399999999
$ws.Range("N29").Value = -1460.75
$ws.Range("H68").Value = 626.7778
$ws.Range("I68").Value = 691.8333
$ws.Range("J68").Value = 496.66666
$ws.Range("K68").Value = 2075.4999
$ws.Range("L68").Value = 1489.99998
$ws.Range("M68").Value = -1264.4999
$ws.Range("N68").Value = -3111.99998
$ws.Range("H71").Value = 626.7778
$ws.Range("I71").Value = 691.8333
$ws.Range("J71").Value = 496.66666
$ws.Range("K71").Value = 6226.4997
$ws.Range("L71").Value = 4469.99994
$ws.Range("M71").Value = -2170.4997
$ws.Range("N71").Value = -12581.99994
$ws.Range("H131").Value = 2013
$ws.Range("J131").Value = 2766.2222
$ws.Range("L131").Value = 8298.6666
$ws.Range("N131").Value = -18378.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 1500
$ws.Range("J23").Value = 1500
$ws.Range("L23").Value = 1500
$ws.Range("N23").Value = -1946
$ws.Range("H54").Value = 49899
$ws.Range("J54").Value = 49899
$ws.Range("L54").Value = 49899
$ws.Range("N54").Value = -50679
$ws.Range("H70").Value = 5747.5
$ws.Range("I70").Value = 5830
$ws.Range("J70").Value = 5500
$ws.Range("K70").Value = 5830
$ws.Range("L70").Value = 5500
$ws.Range("M70").Value = -5560
$ws.Range("N70").Value = -6040
$ws.Range("H73").Value = 5747.5
$ws.Range("I73").Value = 5830
$ws.Range("J73").Value = 5500
$ws.Range("K73").Value = 5830
$ws.Range("L73").Value = 5500
$ws.Range("M73").Value = -4894
$ws.Range("N73").Value = -7372
$ws.Range("H80").Value = 8994.25
$ws.Range("I80").Value = 8996.5
$ws.Range("J80").Value = 8992
$ws.Range("K80").Value = 8996.5
$ws.Range("L80").Value = 8992
$ws.Range("M80").Value = -7998.5
$ws.Range("N80").Value = -10988
$ws.Range("H83").Value = 8994.25
$ws.Range("I83").Value = 8996.5
$ws.Range("J83").Value = 8992
$ws.Range("K83").Value = 44982.5
$ws.Range("L83").Value = 44960
$ws.Range("M83").Value = -39990.5
$ws.Range("N83").Value = -54944
$ws.Range("H132").Value = 9475.799999999999
$ws.Range("I132").Value = 9475.799999999999
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 28427.4
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -25897.4
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3400.6
$ws.Range("I40").Value = 3666
$ws.Range("K40").Value = 3666
$ws.Range("M40").Value = -3530
$ws.Range("H46").Value = 6518.4
$ws.Range("I46").Value = 4036.8
$ws.Range("K46").Value = 4036.8
$ws.Range("M46").Value = -3848.8
$ws.Range("H122").Value = 9360.799999999999
$ws.Range("I122").Value = 9590.315000000001
$ws.Range("K122").Value = 28770.945
$ws.Range("M122").Value = -26320.945

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H126").Value = 2247.5
$ws.Range("I126").Value = 2247.5
$ws.Range("K126").Value = 6742.5
$ws.Range("M126").Value = -4272.5
$ws.Range("H132").Value = 6228.467
$ws.Range("I132").Value = 5365.875
$ws.Range("K132").Value = 16097.625
$ws.Range("M132").Value = -13567.625
